$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 276.44446
$ws.Range("I2").Value2 = 298.5
$ws.Range("K2").Value2 = 298.5
$ws.Range("M2").Value2 = -185.5
$ws.Range("H51").Value2 = 24999.25
$ws.Range("I51").Value2 = 0
$ws.Range("K51").Value2 = 0
$ws.Range("M51").ClearContents()
$ws.Range("H74").Value2 = 14995.625
$ws.Range("I74").Value2 = 15328.667
$ws.Range("K74").Value2 = 15328.667
$ws.Range("M74").Value2 = -14392.667
$ws.Range("H77").Value2 = 14995.625
$ws.Range("I77").Value2 = 15328.667
$ws.Range("K77").Value2 = 76643.33499999999
$ws.Range("M77").Value2 = -71963.33499999999
$ws.Range("H111").Value2 = 4407.0557
$ws.Range("I111").Value2 = 7537.6665
$ws.Range("J111").Value2 = 1276.4445
$ws.Range("K111").Value2 = 22612.9995
$ws.Range("L111").Value2 = 3829.3335
$ws.Range("M111").Value2 = -19545.9995
$ws.Range("N111").Value2 = -9963.333500000001
$ws.Range("H132").Value2 = 4085.4595
$ws.Range("I132").Value2 = 4705.4517
$ws.Range("J132").Value2 = 882.1667
$ws.Range("K132").Value2 = 14116.3551
$ws.Range("L132").Value2 = 2646.5001
$ws.Range("M132").Value2 = -11586.3551
$ws.Range("N132").Value2 = -7706.5001
$ws.Range("H138").Value2 = 2648.9285
$ws.Range("I138").Value2 = 887.3333
$ws.Range("J138").Value2 = 5819.8
$ws.Range("K138").Value2 = 2661.9999
$ws.Range("L138").Value2 = 17459.4
$ws.Range("M138").Value2 = 2478.0001
$ws.Range("N138").Value2 = -27739.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 76.666664
$ws.Range("I5").Value2 = 76.666664
$ws.Range("K5").Value2 = 76.666664
$ws.Range("M5").Value2 = 35.333336
$ws.Range("H32").Value2 = 2568608
$ws.Range("I32").Value2 = 1238117.2
$ws.Range("K32").Value2 = 1238117.2
$ws.Range("M32").Value2 = -1237830.2
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("M34").ClearContents()
$ws.Range("H45").Value2 = 25356.111
$ws.Range("I45").Value2 = 25356.111
$ws.Range("K45").Value2 = 25356.111
$ws.Range("M45").Value2 = -24979.111
$ws.Range("H132").Value2 = 1933.7241
$ws.Range("I132").Value2 = 1639.909
$ws.Range("K132").Value2 = 4919.727000000001
$ws.Range("M132").Value2 = -2389.727000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 76.666664
$ws.Range("I4").Value2 = 76.666664
$ws.Range("K4").Value2 = 76.666664
$ws.Range("M4").Value2 = 38.333336
$ws.Range("H86").Value2 = 4350.3335
$ws.Range("I86").Value2 = 4316.4375
$ws.Range("J86").Value2 = 4458.8
$ws.Range("K86").Value2 = 4316.4375
$ws.Range("L86").Value2 = 4458.8
$ws.Range("M86").Value2 = -3193.4375
$ws.Range("N86").Value2 = -6704.8
$ws.Range("H89").Value2 = 4350.3335
$ws.Range("I89").Value2 = 4316.4375
$ws.Range("J89").Value2 = 4458.8
$ws.Range("K89").Value2 = 21582.1875
$ws.Range("L89").Value2 = 22294
$ws.Range("M89").Value2 = -15966.1875
$ws.Range("N89").Value2 = -33526
$ws.Range("H107").Value2 = 2653842.2
$ws.Range("I107").Value2 = 3345720.5
$ws.Range("K107").Value2 = 3345720.5
$ws.Range("M107").Value2 = -3343800.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1204.3334
$ws.Range("I16").Value2 = 1145.2
$ws.Range("K16").Value2 = 1145.2
$ws.Range("M16").Value2 = -858.2
$ws.Range("H88").Value2 = 20000
$ws.Range("J88").Value2 = 20000
$ws.Range("L88").Value2 = 20000
$ws.Range("N88").Value2 = -20812
$ws.Range("H91").Value2 = 20000
$ws.Range("J91").Value2 = 20000
$ws.Range("L91").Value2 = 20000
$ws.Range("N91").Value2 = -22808
$ws.Range("H94").Value2 = 541.4
$ws.Range("J94").Value2 = 602
$ws.Range("L94").Value2 = 602
$ws.Range("N94").Value2 = -1504
$ws.Range("H107").Value2 = 2381702
$ws.Range("I107").Value2 = 3571968.2
$ws.Range("J107").Value2 = 1169
$ws.Range("K107").Value2 = 3571968.2
$ws.Range("L107").Value2 = 1169
$ws.Range("M107").Value2 = -3570048.2
$ws.Range("N107").Value2 = -5009
$ws.Range("H113").Value2 = 1204.3334
$ws.Range("I113").Value2 = 1145.2
$ws.Range("K113").Value2 = 1145.2
$ws.Range("M113").Value2 = 1024.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 47.25
$ws.Range("J12").Value2 = 62.666668
$ws.Range("L12").Value2 = 188.000004
$ws.Range("N12").Value2 = -534.000004
$ws.Range("H86").Value2 = 1264
$ws.Range("I86").Value2 = 349.5
$ws.Range("K86").Value2 = 1048.5
$ws.Range("M86").Value2 = 137.5
$ws.Range("H89").Value2 = 1264
$ws.Range("I89").Value2 = 349.5
$ws.Range("K89").Value2 = 3145.5
$ws.Range("M89").Value2 = 2782.5
$ws.Range("H131").Value2 = 1641148.5
$ws.Range("I131").Value2 = 28804.143
$ws.Range("J131").Value2 = 2030335.1
$ws.Range("K131").Value2 = 86412.429
$ws.Range("L131").Value2 = 6091005.300000001
$ws.Range("M131").Value2 = -81372.429
$ws.Range("N131").Value2 = -6101085.300000001
$ws.Range("H136").Value2 = 9311.286
$ws.Range("I136").Value2 = 1295.25
$ws.Range("J136").Value2 = 19999.334
$ws.Range("K136").Value2 = 3885.75
$ws.Range("L136").Value2 = 59998.00199999999
$ws.Range("M136").Value2 = 1214.25
$ws.Range("N136").Value2 = -70198.00199999999
$ws.Range("H137").Value2 = 1754
$ws.Range("I137").Value2 = 1704.8
$ws.Range("J137").Value2 = 2000
$ws.Range("K137").Value2 = 5114.4
$ws.Range("L137").Value2 = 6000
$ws.Range("M137").Value2 = -14.39999999999964
$ws.Range("N137").Value2 = -16200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value2 = 100000
$ws.Range("I12").Value2 = 100000
$ws.Range("K12").Value2 = 100000
$ws.Range("M12").Value2 = -99860
$ws.Range("H13").Value2 = 2276.25
$ws.Range("I13").Value2 = 25
$ws.Range("J13").Value2 = 3026.6667
$ws.Range("K13").Value2 = 25
$ws.Range("L13").Value2 = 3026.6667
$ws.Range("M13").Value2 = 114
$ws.Range("N13").Value2 = -3304.6667
$ws.Range("H17").Value2 = 1152.25
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 1152.25
$ws.Range("K17").Value2 = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value2 = 1152.25
$ws.Range("N17").Value2 = -1488.25
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("M19").ClearContents()
$ws.Range("H20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value2 = 0
$ws.Range("H22").Value2 = 755.5
$ws.Range("J22").Value2 = 1111
$ws.Range("L22").Value2 = 1111
$ws.Range("N22").Value2 = -2169
$ws.Range("H23").Value2 = 390
$ws.Range("J23").Value2 = 390
$ws.Range("L23").Value2 = 390
$ws.Range("N23").Value2 = -836
$ws.Range("H24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value2 = 0
$ws.Range("H25").Value2 = 4250
$ws.Range("I25").Value2 = 1500
$ws.Range("J25").Value2 = 7000
$ws.Range("K25").Value2 = 1500
$ws.Range("L25").Value2 = 7000
$ws.Range("M25").Value2 = -971
$ws.Range("N25").Value2 = -8058
$ws.Range("H80").Value2 = 90913290
$ws.Range("J80").Value2 = 4577.4443
$ws.Range("L80").Value2 = 4577.4443
$ws.Range("N80").Value2 = -6573.4443
$ws.Range("H83").Value2 = 90913290
$ws.Range("J83").Value2 = 4577.4443
$ws.Range("L83").Value2 = 22887.2215
$ws.Range("N83").Value2 = -32871.2215
$ws.Range("H107").Value2 = 2645.1538
$ws.Range("I107").Value2 = 2017.2858
$ws.Range("K107").Value2 = 2017.2858
$ws.Range("M107").Value2 = -97.28580000000011
$ws.Range("H113").Value2 = 2123.3157
$ws.Range("I113").Value2 = 2096
$ws.Range("K113").Value2 = 2096
$ws.Range("M113").Value2 = 74
$ws.Range("H132").Value2 = 2530.4
$ws.Range("I132").Value2 = 2243.4285
$ws.Range("K132").Value2 = 6730.2855
$ws.Range("M132").Value2 = -4200.2855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 478.58823
$ws.Range("I55").Value2 = 270.77777
$ws.Range("K55").Value2 = 270.77777
$ws.Range("M55").Value2 = -97.77776999999998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value2 = 21751.4
$ws.Range("J43").Value2 = 29910
$ws.Range("L43").Value2 = 29910
$ws.Range("N43").Value2 = -30208
$ws.Range("H58").Value2 = 9085
$ws.Range("I58").Value2 = 9085
$ws.Range("K58").Value2 = 9085
$ws.Range("M58").Value2 = -8777
$ws.Range("H107").Value2 = 725.4666999999999
$ws.Range("I107").Value2 = 700.38464
$ws.Range("K107").Value2 = 2101.15392
$ws.Range("M107").Value2 = -181.1539199999997
$ws.Range("H113").Value2 = 451.15384
$ws.Range("I113").Value2 = 447.6
$ws.Range("K113").Value2 = 1342.8
$ws.Range("M113").Value2 = 827.1999999999998
